{"js": "// Helper: build a minimal single-part OOXML \"Flat OPC\" package wrapping a\n// <w:document><w:body>...</w:body></w:document> fragment, suitable for\n// Range/Paragraph.insertOoxml().\nfunction wrapBody(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyInnerXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Edit 1: \"features \" paragraph becomes the \"App name can be ...\" ---\n// paragraph (4 runs) followed by a new \"Features \" paragraph (2 runs).\nlet featuresParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"features \") {\n    featuresParagraph = p;\n    break;\n  }\n}\n\nif (featuresParagraph) {\n  // Insert the new \"Features \" paragraph right after the current one first\n  // (its own text is a placeholder; it gets its final run split below),\n  // so the original paragraph is still intact when we locate it.\n  const newFeaturesParagraph = featuresParagraph.insertParagraph(\"Features \", \"After\");\n  await context.sync();\n\n  // Replace the original paragraph's single run with the four runs that\n  // make up the app-name suggestion.\n  const appNameOoxml = wrapBody(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">App name can be </w:t></w:r>' +\n      \"<w:r><w:t>CTA:</w:t></w:r>\" +\n      '<w:r><w:t xml:space=\"preserve\"> Chandigarh Tourism App because in my view the name should be such that it is self explanatory </w:t></w:r>' +\n      \"<w:r><w:t>that what the app is for</w:t></w:r>\" +\n    \"</w:p>\"\n  );\n  featuresParagraph.insertOoxml(appNameOoxml, \"Replace\");\n  await context.sync();\n\n  // Split the new paragraph's text into two runs: \"Features\" + \" \".\n  const featuresOoxml = wrapBody(\n    \"<w:p>\" +\n      \"<w:r><w:t>Features</w:t></w:r>\" +\n      '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    \"</w:p>\"\n  );\n  newFeaturesParagraph.insertOoxml(featuresOoxml, \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: \"These are some feature please add more if you want .\" ---\n// loses the space before the period and is split into two runs.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet wantParagraph = null;\nfor (const p of paragraphs2.items) {\n  if (p.text === \"These are some feature please add more if you want .\") {\n    wantParagraph = p;\n    break;\n  }\n}\n\nif (wantParagraph) {\n  const wantOoxml = wrapBody(\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/></w:pPr>' +\n      '<w:r><w:t xml:space=\"preserve\">These are some feature please add more if you </w:t></w:r>' +\n      \"<w:r><w:t>want.</w:t></w:r>\" +\n    \"</w:p>\"\n  );\n  wantParagraph.insertOoxml(wantOoxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Helper: wrap a <w:body> fragment as a minimal Flat-OPC \"WordOpenXML\"\n# package suitable for Range.InsertXML().\nfunction Wrap-Body([string]$bodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n        '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n# --- Edit 1: the \"features \" paragraph becomes a 4-run \"App name can be\n#     CTA: Chandigarh Tourism App ...\" paragraph, immediately followed by a\n#     new two-run \"Features \" paragraph. ---\n$featuresParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"features `r\") {\n        $featuresParagraph = $p\n        break\n    }\n}\n\nif ($featuresParagraph -ne $null) {\n    # Split off a new (empty) paragraph right after the current one.\n    $featuresParagraph.Range.InsertParagraphAfter()\n\n    # Re-locate both paragraphs now that the paragraph collection shifted.\n    $appNameParagraph = $null\n    $newFeaturesParagraph = $null\n    $matchedFirst = $false\n    foreach ($p in $d.Paragraphs) {\n        if (-not $matchedFirst -and $p.Range.Text -eq \"features `r\") {\n            $appNameParagraph = $p\n            $matchedFirst = $true\n            continue\n        }\n        if ($matchedFirst) {\n            $newFeaturesParagraph = $p\n            break\n        }\n    }\n\n    $appNameXml = Wrap-Body(\n        '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">App name can be </w:t></w:r>' +\n        '<w:r><w:t>CTA:</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Chandigarh Tourism App because in my view the name should be such that it is self explanatory </w:t></w:r>' +\n        '<w:r><w:t>that what the app is for</w:t></w:r>' +\n        '</w:p>'\n    )\n    $appNameParagraph.Range.InsertXML($appNameXml)\n\n    $featuresXml = Wrap-Body(\n        '<w:p>' +\n        '<w:r><w:t>Features</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '</w:p>'\n    )\n    $newFeaturesParagraph.Range.InsertXML($featuresXml)\n}\n\n# --- Edit 2: \"These are some feature please add more if you want .\" loses\n#     the space before the period and is split into two runs. ---\n$wantParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"These are some feature please add more if you want .`r\") {\n        $wantParagraph = $p\n        break\n    }\n}\n\nif ($wantParagraph -ne $null) {\n    $wantXml = Wrap-Body(\n        '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/></w:pPr>' +\n        '<w:r><w:t xml:space=\"preserve\">These are some feature please add more if you </w:t></w:r>' +\n        '<w:r><w:t>want.</w:t></w:r>' +\n        '</w:p>'\n    )\n    $wantParagraph.Range.InsertXML($wantXml)\n}\n"}
